# Update the "as_of_utc" timestamp column (AA) for rows 2-26 on both the
# "Главные" (sheet 2) and "Линейные" (sheet 3) worksheets, changing the
# previously captured timestamp "2025-11-10 03:06:39" to the refreshed
# publish timestamp "2025-11-10 07:06:41".

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-10 07:06:41"
$col = 27  # column AA

for ($sheetIndex = 2; $sheetIndex -le 3; $sheetIndex++) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, $col).Value = $newTimestamp
    }
}
